$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja 1")

$ws.Range("B17").Value = 'El costo del dispositivo debe ser menor o igual a 50 U\$D.'
$ws.Range("B17").Select() | Out-Null
